$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching style of existing headers (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-32
$data = @(
    @(2, 5),
    @(8, 8),
    @(7, 8),
    @(8, 9),
    @(1, 4),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(9, 9),
    @(4, 7),
    @(2, 5),
    @(4, 6),
    @(7, 7),
    @(7, 7),
    @(8, 9),
    @(7, 8),
    @(11, 12),
    @(7, 8),
    @(1, 3),
    @(6, 8),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 3),
    @(3, 4),
    @(1, 3),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
